{"js": "// Replace the division problems in the worksheet table with the new values,\n// matching the unified diff. The replacement list is in document order, so\n// for a source value that appears more than once (e.g. \"89\u00f74=\"), each\n// successive entry finds the NEXT remaining occurrence -- once one instance\n// is rewritten to its new text it naturally drops out of later searches for\n// the same old text.\nconst replacements = [\n  [\"43\u00f78=\", \"86\u00f72=\"],\n  [\"64\u00f78=\", \"17\u00f77=\"],\n  [\"89\u00f78=\", \"88\u00f74=\"],\n  [\"72\u00f75=\", \"45\u00f74=\"],\n  [\"54\u00f78=\", \"28\u00f75=\"],\n  [\"26\u00f73=\", \"52\u00f78=\"],\n  [\"65\u00f72=\", \"95\u00f72=\"],\n  [\"93\u00f76=\", \"73\u00f73=\"],\n  [\"89\u00f74=\", \"46\u00f76=\"],\n  [\"25\u00f79=\", \"56\u00f76=\"],\n  [\"75\u00f72=\", \"22\u00f78=\"],\n  [\"37\u00f75=\", \"36\u00f75=\"],\n  [\"37\u00f78=\", \"34\u00f75=\"],\n  [\"67\u00f72=\", \"59\u00f73=\"],\n  [\"15\u00f75=\", \"49\u00f74=\"],\n  [\"85\u00f74=\", \"55\u00f72=\"],\n  [\"24\u00f76=\", \"92\u00f76=\"],\n  [\"24\u00f79=\", \"71\u00f79=\"],\n  [\"26\u00f76=\", \"38\u00f74=\"],\n  [\"28\u00f74=\", \"62\u00f77=\"],\n  [\"89\u00f74=\", \"98\u00f72=\"],\n  [\"85\u00f78=\", \"48\u00f74=\"],\n  [\"75\u00f75=\", \"20\u00f75=\"],\n  [\"63\u00f77=\", \"33\u00f79=\"],\n  [\"80\u00f72=\", \"96\u00f73=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Could not find text \"${oldText}\" to replace.`);\n  }\n\n  // The first remaining hit is always the right one: earlier replacements in\n  // this list already rewrote any prior occurrences of this same old text,\n  // so they no longer match this search.\n  results.items[0].insertText(newText, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Replace the division problems in the worksheet table with the new values,\n# matching the unified diff.\n#\n# Some \"old\" values repeat (e.g. \"89\u00f74=\" appears twice, mapping to two\n# different replacements), so instead of one global \"replace all\" we issue a\n# fresh Find/Replace (wdReplaceOne) for every pair, in document order. Each\n# call only touches the first remaining match of that old text; once an\n# occurrence is rewritten it no longer matches later searches for the same\n# old text, so subsequent entries correctly land on the next instance.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"43\u00f78=\", \"86\u00f72=\"),\n    @(\"64\u00f78=\", \"17\u00f77=\"),\n    @(\"89\u00f78=\", \"88\u00f74=\"),\n    @(\"72\u00f75=\", \"45\u00f74=\"),\n    @(\"54\u00f78=\", \"28\u00f75=\"),\n    @(\"26\u00f73=\", \"52\u00f78=\"),\n    @(\"65\u00f72=\", \"95\u00f72=\"),\n    @(\"93\u00f76=\", \"73\u00f73=\"),\n    @(\"89\u00f74=\", \"46\u00f76=\"),\n    @(\"25\u00f79=\", \"56\u00f76=\"),\n    @(\"75\u00f72=\", \"22\u00f78=\"),\n    @(\"37\u00f75=\", \"36\u00f75=\"),\n    @(\"37\u00f78=\", \"34\u00f75=\"),\n    @(\"67\u00f72=\", \"59\u00f73=\"),\n    @(\"15\u00f75=\", \"49\u00f74=\"),\n    @(\"85\u00f74=\", \"55\u00f72=\"),\n    @(\"24\u00f76=\", \"92\u00f76=\"),\n    @(\"24\u00f79=\", \"71\u00f79=\"),\n    @(\"26\u00f76=\", \"38\u00f74=\"),\n    @(\"28\u00f74=\", \"62\u00f77=\"),\n    @(\"89\u00f74=\", \"98\u00f72=\"),\n    @(\"85\u00f78=\", \"48\u00f74=\"),\n    @(\"75\u00f75=\", \"20\u00f75=\"),\n    @(\"63\u00f77=\", \"33\u00f79=\"),\n    @(\"80\u00f72=\", \"96\u00f73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n\n    $found = $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, \"wdFindContinue\", $false, $find.Replacement.Text, \"wdReplaceOne\")\n\n    if (-not $found) {\n        Write-Output \"WARNING: could not find '$oldText' to replace with '$newText'\"\n    }\n}\n"}
